# "writing 5 most recent values is now working"
# Fill in the two patient rows (row 2 = mayar, row 3 = Shehab) with the
# 5 most recent heart-rate / temperature readings plus blood type, sex,
# age and the date-created / date-modified timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---- Row 2 : patient "mayar" -------------------------------------------------
$ws.Cells.Item(2, 1).Value = 1            # A2  ID
$ws.Cells.Item(2, 2).Value = "mayar"      # B2  Name

$ws.Cells.Item(2, 3).Value = 107          # C2  Heart Rate1
$ws.Cells.Item(2, 3).Style = "Normal"
$ws.Cells.Item(2, 4).Value = 80           # D2  Heart Rate2
$ws.Cells.Item(2, 4).Style = "Normal"
$ws.Cells.Item(2, 5).Value = 89           # E2  Heart Rate3
$ws.Cells.Item(2, 5).Style = "Normal"
$ws.Cells.Item(2, 6).Value = 105          # F2  Heart Rate4
$ws.Cells.Item(2, 6).Style = "Normal"
$ws.Cells.Item(2, 7).Value = 75           # G2  Heart Rate5
$ws.Cells.Item(2, 7).Style = "Normal"

$ws.Cells.Item(2, 8).Value = 36.16999816894531    # H2  Temp1
$ws.Cells.Item(2, 8).Style = "Normal"
$ws.Cells.Item(2, 9).Value = 36.70000076293945    # I2  Temp2
$ws.Cells.Item(2, 9).Style = "Normal"
$ws.Cells.Item(2, 10).Value = 35.709999084472656  # J2  Temp3
$ws.Cells.Item(2, 10).Style = "Normal"
$ws.Cells.Item(2, 11).Value = 35.70000076293945   # K2  Temp4
$ws.Cells.Item(2, 11).Style = "Normal"
$ws.Cells.Item(2, 12).Value = 36.900001525878906  # L2  Temp5
$ws.Cells.Item(2, 12).Style = "Normal"

$ws.Cells.Item(2, 13).Value = "B-"        # M2  Blood Type
$ws.Cells.Item(2, 14).Value = "female"    # N2  Sex
$ws.Cells.Item(2, 15).Value = 12          # O2  Age

$ws.Cells.Item(2, 16).Value = 42501.207151319446  # P2  Date Created
$ws.Cells.Item(2, 16).NumberFormat = "d/m/yy hh:mm:ss"
$ws.Cells.Item(2, 17).Value = 42501.22014013889   # Q2  Date Modified
$ws.Cells.Item(2, 17).NumberFormat = "d/m/yy hh:mm:ss"

# ---- Row 3 : patient "Shehab" ------------------------------------------------
$ws.Cells.Item(3, 1).Value = 2            # A3  ID
$ws.Cells.Item(3, 2).Value = "Shehab"     # B3  Name

$ws.Cells.Item(3, 3).Value = 94           # C3  Heart Rate1
$ws.Cells.Item(3, 3).Style = "Normal"
$ws.Cells.Item(3, 4).Value = 69           # D3  Heart Rate2
$ws.Cells.Item(3, 4).Style = "Normal"
$ws.Cells.Item(3, 5).Value = 107          # E3  Heart Rate3
$ws.Cells.Item(3, 5).Style = "Normal"
$ws.Cells.Item(3, 6).Value = 108          # F3  Heart Rate4
$ws.Cells.Item(3, 6).Style = "Normal"
$ws.Cells.Item(3, 7).Value = 71           # G3  Heart Rate5
$ws.Cells.Item(3, 7).Style = "Normal"

$ws.Cells.Item(3, 8).Value = 35.86000061035156    # H3  Temp1
$ws.Cells.Item(3, 8).Style = "Normal"
$ws.Cells.Item(3, 9).Value = 35.77000045776367    # I3  Temp2
$ws.Cells.Item(3, 9).Style = "Normal"
$ws.Cells.Item(3, 10).Value = 37.09000015258789   # J3  Temp3
$ws.Cells.Item(3, 10).Style = "Normal"
$ws.Cells.Item(3, 11).Value = 37.029998779296875  # K3  Temp4
$ws.Cells.Item(3, 11).Style = "Normal"
$ws.Cells.Item(3, 12).Value = 36.0099983215332    # L3  Temp5
$ws.Cells.Item(3, 12).Style = "Normal"

$ws.Cells.Item(3, 13).Value = "O+"        # M3  Blood Type
$ws.Cells.Item(3, 14).Value = "male"      # N3  Sex
$ws.Cells.Item(3, 15).Value = 25          # O3  Age

$ws.Cells.Item(3, 16).Value = 42502.77497329861   # P3  Date Created
$ws.Cells.Item(3, 16).NumberFormat = "d/m/yy hh:mm:ss"

# Row 3 no longer has a "Date Modified" value - remove it entirely.
$ws.Cells.Item(3, 17).Style = "Normal"
$ws.Cells.Item(3, 17).ClearContents()

# ---- Selection: column R ("Alarm") is now the active selection --------------
$ws.Range("R1:R1048576").Select()
